# Weekly Fruta/Hortaliza update:
# Insert two new rows of data (2022-02-25) above the existing row 182,
# pushing the old rows 182-203 down to 184-205.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 182 (existing row 182 and everything below
# shifts down by 2).
$ws.Rows.Item(182).Insert()
$ws.Rows.Item(182).Insert()

# New row 182: Red Globe
$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = "2022-02-25"
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = "Fruta"
$ws.Range("G182").Value = 100109
$ws.Range("H182").Value = "Uva"
$ws.Range("I182").Value = 100109001
$ws.Range("J182").Value = "Uva"
$ws.Range("K182").Value = "Red Globe"
$ws.Range("L182").Value = "Primera"
$ws.Range("M182").Value = 300
$ws.Range("N182").Value = 17000
$ws.Range("O182").Value = 18000
$ws.Range("P182").Value = 17500
$ws.Range("Q182").Value = "`$/caja 20 kilos"
$ws.Range("R182").Value = "Región de O'Higgins"
$ws.Range("S182").Value = 875
$ws.Range("T182").Value = 20

# New row 183: Superior Seedless
$ws.Range("A183").Value = 4
$ws.Range("B183").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C183").Value = "Los Lagos"
$ws.Range("D183").Value = "2022-02-25"
$ws.Range("E183").Value = 10
$ws.Range("F183").Value = "Fruta"
$ws.Range("G183").Value = 100109
$ws.Range("H183").Value = "Uva"
$ws.Range("I183").Value = 100109001
$ws.Range("J183").Value = "Uva"
$ws.Range("K183").Value = "Superior Seedless"
$ws.Range("L183").Value = "Primera"
$ws.Range("M183").Value = 200
$ws.Range("N183").Value = 17000
$ws.Range("O183").Value = 18000
$ws.Range("P183").Value = 17500
$ws.Range("Q183").Value = "`$/caja 20 kilos"
$ws.Range("R183").Value = "Región de O'Higgins"
$ws.Range("S183").Value = 875
$ws.Range("T183").Value = 20
